$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column B values for rows 2-118 (Epoch Accuracy values), per the
# latest "Froze Encoder 12" run.
$bValues = @(
    0.796875, 0.703125, 0.65625, 0.640625, 0.546875, 0.53125, 0.5, 0.484375,
    0.515625, 0.5, 0.453125, 0.5625, 0.5625, 0.5, 0.515625, 0.40625,
    0.421875, 0.421875, 0.421875, 0.421875, 0.421875, 0.421875, 0.421875, 0.421875,
    0.421875, 0.421875, 0.421875, 0.421875, 0.421875, 0.421875, 0.4375, 0.4375,
    0.453125, 0.453125, 0.453125, 0.453125, 0.46875, 0.46875, 0.46875, 0.46875,
    0.453125, 0.453125, 0.453125, 0.453125, 0.453125, 0.453125, 0.453125, 0.453125,
    0.453125, 0.453125, 0.453125, 0.453125, 0.453125, 0.453125, 0.453125, 0.453125,
    0.453125, 0.453125, 0.453125, 0.453125, 0.453125, 0.453125, 0.453125, 0.453125,
    0.453125, 0.453125, 0.453125, 0.453125, 0.453125, 0.453125, 0.453125, 0.453125,
    0.453125, 0.453125, 0.453125, 0.453125, 0.453125, 0.4375, 0.4375, 0.4375,
    0.4375, 0.4375, 0.453125, 0.453125, 0.453125, 0.453125, 0.453125, 0.46875,
    0.46875, 0.46875, 0.46875, 0.46875, 0.46875, 0.46875, 0.46875, 0.46875,
    0.46875, 0.46875, 0.46875, 0.46875, 0.46875, 0.484375, 0.34375, 0.515625,
    0.28125, 0.328125, 0.359375, 0.46875, 0.4375, 0.453125, 0.546875, 0.46875,
    0.296875, 0.375, 0.40625, 0.375, (23.0/61.0)
)

for ($i = 0; $i -lt $bValues.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 2).Value = $bValues[$i]
}

# The DisplayOutputs object repr in column A (rows 102-118) carries a new
# memory address after the re-run.
$newRepr = "<__main__.DisplayOutputs object at 0x7f94340732b0>"
for ($row = 102; $row -le 118; $row++) {
    $ws.Cells.Item($row, 1).Value = $newRepr
}

# Selection changed from A2:B116 to the full-sheet range A1:XFD1048576.
$ws.Range("A1:XFD1048576").Select()
